$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.136329588014981
$ws.Range("C2").Value = 0.0284644194756554
$ws.Range("D2").Value = 0.0262172284644195
$ws.Range("E2").Value = 0.964794007490637
$ws.Range("F2").Value = 0.0149812734082397
$ws.Range("G2").Value = 0.941573033707865
$ws.Range("H2").Value = 0.0239700374531835
$ws.Range("I2").Value = 0.0194756554307116
$ws.Range("J2").Value = 0.0434456928838951
$ws.Range("K2").Value = 0.946067415730337
$ws.Range("L2").Value = 0.0441947565543071
$ws.Range("M2").Value = 0.115355805243446
$ws.Range("O2").Value = 0.0471910112359551
$ws.Range("P2").Value = 0.0202247191011236
$ws.Range("Q2").Value = 0.00449438202247191
$ws.Range("R2").Value = 0.00149812734082397
$ws.Range("S2").Value = 0.953558052434457
$ws.Range("T2").Value = 0.00149812734082397
$ws.Range("U2").Value = 0.976779026217228
$ws.Range("V2").Value = 0.953558052434457
$ws.Range("W2").Value = 0.0262172284644195
$ws.Range("X2").Value = 0.0247191011235955

$ws.Range("B3").Value = 0.811985018726592
$ws.Range("C3").Value = 0.0217228464419476
$ws.Range("D3").Value = 0.0052434456928839
$ws.Range("E3").Value = 0.00449438202247191
$ws.Range("F3").Value = 0.943820224719101
$ws.Range("G3").Value = 0.050187265917603
$ws.Range("H3").Value = 0.161797752808989
$ws.Range("I3").Value = 0.856179775280899
$ws.Range("J3").Value = 0.0397003745318352
$ws.Range("K3").Value = 0.00599250936329588
$ws.Range("M3").Value = 0.832958801498127
$ws.Range("N3").Value = 0.0194756554307116
$ws.Range("O3").Value = 0.00299625468164794
$ws.Range("P3").Value = 0.00224719101123596
$ws.Range("Q3").Value = 0.991760299625468
$ws.Range("R3").Value = 0.0142322097378277
$ws.Range("S3").Value = 0.00149812734082397
$ws.Range("T3").Value = 0.9812734082397
$ws.Range("U3").Value = 0.00973782771535581
$ws.Range("V3").Value = 0.0142322097378277
$ws.Range("W3").Value = 0.0127340823970037
$ws.Range("X3").Value = 0.0119850187265918

$ws.Range("B4").Value = 0.0456928838951311
$ws.Range("C4").Value = 0.00149812734082397
$ws.Range("D4").Value = 0.953558052434457
$ws.Range("E4").Value = 0.0307116104868914
$ws.Range("F4").Value = 0.00299625468164794
$ws.Range("G4").Value = 0.000749063670411985
$ws.Range("H4").Value = 0.0239700374531835
$ws.Range("I4").Value = 0.0142322097378277
$ws.Range("J4").Value = 0.91310861423221
$ws.Range("K4").Value = 0.0419475655430712
$ws.Range("L4").Value = 0.954307116104869
$ws.Range("M4").Value = 0.00299625468164794
$ws.Range("N4").Value = 0.0209737827715356
$ws.Range("P4").Value = 0.188014981273408
$ws.Range("Q4").Value = 0.00149812734082397
$ws.Range("R4").Value = 0.000749063670411985
$ws.Range("S4").Value = 0.0434456928838951
$ws.Range("T4").Value = 0.00749063670411985
$ws.Range("U4").Value = 0.00898876404494382
$ws.Range("V4").Value = 0.0292134831460674
$ws.Range("W4").Value = 0.959550561797753
$ws.Range("X4").Value = 0.961048689138577

$ws.Range("B5").Value = 0.00599250936329588
$ws.Range("C5").Value = 0.948314606741573
$ws.Range("D5").Value = 0.0149812734082397
$ws.Range("F5").Value = 0.0382022471910112
$ws.Range("G5").Value = 0.00749063670411985
$ws.Range("H5").Value = 0.790262172284644
$ws.Range("I5").Value = 0.110112359550562
$ws.Range("J5").Value = 0.00374531835205993
$ws.Range("K5").Value = 0.00599250936329588
$ws.Range("L5").Value = 0.00149812734082397
$ws.Range("M5").Value = 0.048689138576779
$ws.Range("N5").Value = 0.959550561797753
$ws.Range("O5").Value = 0.949812734082397
$ws.Range("P5").Value = 0.789513108614232
$ws.Range("Q5").Value = 0.00224719101123596
$ws.Range("R5").Value = 0.983520599250936
$ws.Range("S5").Value = 0.00149812734082397
$ws.Range("T5").Value = 0.00973782771535581
$ws.Range("U5").Value = 0.00449438202247191
$ws.Range("V5").Value = 0.00299625468164794
$ws.Range("W5").Value = 0.000749063670411985
$ws.Range("X5").Value = 0.00224719101123596
